$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column S (old S/"!EndOfField" marker column
# shifts right to become column T), making room for the new "max_health"
# property column between "attack_interval" (renamed from "max_life", now
# in column R) and the end-of-table marker.
$ws.Range("S1").EntireColumn.Insert()

# --- Row 2 (field/column names) ---
# Column R used to hold "max_life"; rename it to "attack_interval".
$ws.Range("R2").Value2 = "attack_interval"
# New column S holds the new "max_health" field name.
$ws.Range("S2").Value2 = "max_health"
$ws.Range("S2").Style = $ws.Range("R2").Style

# --- Row 3 (type markers) ---
# New column S is a Float field, same as its neighbours.
$ws.Range("S3").Value2 = "!Float"
$ws.Range("S3").Style = $ws.Range("R3").Style

# --- Row 4 (None entity - all blank) ---
$ws.Range("S4").Style = $ws.Range("R4").Style

# --- Data rows: attack_power (Q), attack_interval (R), max_health (S) ---
# Bullet (row 5)
$ws.Range("Q5").Value2 = 5
$ws.Range("R5").Value2 = ""
$ws.Range("S5").Value2 = 1
$ws.Range("S5").Style = $ws.Range("R5").Style

# Hunter (row 6)
$ws.Range("Q6").Value2 = 5
$ws.Range("R6").Value2 = 0.3
$ws.Range("S6").Value2 = 50
$ws.Range("S6").Style = $ws.Range("R6").Style

# Bot_X (row 7)
$ws.Range("Q7").Value2 = 5
$ws.Range("R7").Value2 = 0.3
$ws.Range("S7").Value2 = 50
$ws.Range("S7").Style = $ws.Range("R7").Style

# Player (row 8)
$ws.Range("Q8").Value2 = 5
$ws.Range("R8").Value2 = 0.3
$ws.Range("S8").Value2 = 100
$ws.Range("S8").Style = $ws.Range("R8").Style

# --- Row 9 (end-of-table marker row - all blank) ---
$ws.Range("S9").Style = $ws.Range("R9").Style

# --- Column widths ---
# "attack_interval" is wider than "max_life" was, so widen column R and
# drop its best-fit autosize flag; give the new column S its own best-fit
# width.
$ws.Columns("R").ColumnWidth = 12.296875
$ws.Columns("S").ColumnWidth = 10.3984375

# --- Merged header cell above attack_power/attack_interval/max_health ---
$ws.Range("Q1:R1").UnMerge()
$ws.Range("Q1:S1").Merge()

# --- View state ---
$ws.Application.ActiveWindow.ScrollColumn = 6
$ws.Range("O7").Select()
